# Rename the two SQL-code-gen worksheets and make the (renamed) "add
# column" sheet the active tab, matching "rename sheets in sql code gen".
$wb = $excel.ActiveWorkbook

$wsAdd  = $wb.Worksheets.Item("SQL-Code - Add new custom field")
$wsDrop = $wb.Worksheets.Item("SQL-Code - Drop custom field")

$wsAdd.Name  = "CodeGen-AddColumn"
$wsDrop.Name = "CodeGen-DeleteColumn"

$wsAdd.Activate()
